# Apply updated cryptocurrency price/volume data to Sheet1
# (matches the "Updated cryptos list" GitHub Actions commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'" + "54.205.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + "  -3.21%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'" + "2.286.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + "  -3.25%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'" + "  -0.05%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'" + "492.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + "  -1.86%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'" + "128.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + "  -2.13%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'" + "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + "  -0.20%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'" + "0.527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + "  -3.81%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'" + "2.292.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + "  -2.78%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'" + "0.0938"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + "  -3.49%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'" + "0.148"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + "  -1.23%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'" + "4.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + "  +3.21%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'" + "0.318"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + "  -2.84%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'" + "2.692.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + "  -3.23%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "'" + "  +0.00%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'" + "54.190.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + "  -3.12%  "
$ws.Range("E16").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'" + "2.267.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + "  -4.60%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'" + "  -3.00%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'" + "4.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + "  +0.29%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'" + "303.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + "  -0.82%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'" + "6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + "  -0.25%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "'" + "  +0.03%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'" + "64.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'" + "  -1.33%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'" + "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + "  -0.33%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'" + "0.367"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + "  -0.70%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'" + "0.144"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + "  -2.49%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'" + "  -1.17%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'" + "169.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'" + "  -1.46%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.0" + [char]0x2083 + "0700"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + "  -1.98%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'" + "  -1.18%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'" + "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'" + "  -0.09%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'" + "5.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'" + "  +1.13%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'" + "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + "  +0.12%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "'" + "  -1.66%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'" + "  +0.30%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'" + "  -0.93%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'" + "0.854"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + "  +7.68%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = "'" + "  -3.81%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'" + "35.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + "  -0.64%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'" + "  -1.95%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'" + "0.368"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + "  -0.23%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'" + "3.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + "  -0.34%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'" + "123.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + "  -5.73%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'" + "4.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + "  -0.68%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "'" + "  -2.62%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'" + "0.545"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + "  -2.79%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'" + "240.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + "  -1.02%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'" + "0.0475"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'" + "  -0.65%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'" + "  -1.35%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'" + "16.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + "  -2.07%  "
$ws.Range("E51").Style = "Normal"

